# Fill in the missing/blank Column B (most recent period) figures on the
# "CB" cash-flow sheet, correct a handful of mis-entered historical values,
# and widen column B to match column C's width.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CB")

# --- Column B: previously-blank cells now populated with values ---------
$ws.Range("B3").Value  = 645000000
$ws.Range("B4").Value  = -1667000000
$ws.Range("B5").Value  = -1059000000
$ws.Range("B6").Value  = 7341000000
$ws.Range("B7").Value  = 38000000
$ws.Range("B9").Value  = -701000000
$ws.Range("B10").Value = 10874000000
$ws.Range("B12").Value = -38000000
$ws.Range("B13").Value = -6248000000
$ws.Range("B14").Value = 967000000
$ws.Range("B15").Value = -8459000000
$ws.Range("B17").Value = -574000000
$ws.Range("B18").Value = -1401000000
$ws.Range("B19").Value = 86000000
$ws.Range("B20").Value = -2249000000
$ws.Range("B21").Value = 17000000
$ws.Range("B22").Value = 183000000
$ws.Range("B23").Value = 1836000000
$ws.Range("B24").Value = 2019000000
$ws.Range("B26").Value = 6320000000
$ws.Range("B27").Value = -6286000000
$ws.Range("B28").Value = -574000000

# --- Corrections to existing historical values --------------------------
$ws.Range("E11").Value = -3129000000
$ws.Range("F13").Value = -4959000000
$ws.Range("F17").Value = -1753000000
$ws.Range("B25").Value = -1401000000
$ws.Range("F27").Value = -5485000000
$ws.Range("F28").Value = -1753000000

# --- Widen column B to match column C (15.4 -> 16.5) ---------------------
$ws.Columns.Item(2).ColumnWidth = 16.5
